$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.533107333333334
$ws.Range("H2").Value = 22.599322
$ws.Range("I2").Value = 0.4772251808959424
$ws.Range("J2").Value = 0.4772251808959424
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.539707
$ws.Range("N2").Value = 7.619121
$ws.Range("O2").Value = 0.04480768326120512
$ws.Range("P2").Value = 0.04480768326120513
$ws.Range("Q2").Value = 19.131885426218
$ws.Range("R2").Value = 172.186968835962
$ws.Range("S2").Value = 0.0213833547498567
$ws.Range("T2").Value = 0.02138335474985671
$ws.Range("G3").Value = 7.533107333333334
$ws.Range("H3").Value = 22.599322
$ws.Range("I3").Value = 0.4772251808959424
$ws.Range("J3").Value = 0.4772251808959424
$ws.Range("O3").Value = 0.3622039450212636
$ws.Range("P3").Value = 0.3622039450212636
$ws.Range("Q3").Value = 154.6530387807558
$ws.Range("R3").Value = 1391.877349026802
$ws.Range("S3").Value = 0.1728528431839965
$ws.Range("T3").Value = 0.1728528431839965
$ws.Range("G4").Value = 7.533107333333334
$ws.Range("H4").Value = 22.599322
$ws.Range("I4").Value = 0.4772251808959424
$ws.Range("J4").Value = 0.4772251808959424
$ws.Range("O4").Value = 0.5929883717175313
$ws.Range("P4").Value = 0.5929883717175314
$ws.Range("Q4").Value = 253.1928624973558
$ws.Range("R4").Value = 2278.735762476202
$ws.Range("S4").Value = 0.2829889829620892
$ws.Range("T4").Value = 0.2829889829620892
$ws.Range("G5").Value = 5.009378000000001
$ws.Range("I5").Value = 0.3173459790819593
$ws.Range("J5").Value = 0.3173459790819593
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.539707
$ws.Range("N5").Value = 7.619121
$ws.Range("O5").Value = 0.04480768326120512
$ws.Range("P5").Value = 0.04480768326120513
$ws.Range("Q5").Value = 12.722352372246
$ws.Range("R5").Value = 114.501171350214
$ws.Range("S5").Value = 0.01421953811492146
$ws.Range("T5").Value = 0.01421953811492146
$ws.Range("G6").Value = 5.009378000000001
$ws.Range("I6").Value = 0.3173459790819593
$ws.Range("J6").Value = 0.3173459790819593
$ws.Range("O6").Value = 0.3622039450212636
$ws.Range("P6").Value = 0.3622039450212636
$ws.Range("R6").Value = 925.5728695196941
$ws.Range("S6").Value = 0.1149439655601211
$ws.Range("T6").Value = 0.114943965560121
$ws.Range("G7").Value = 5.009378000000001
$ws.Range("I7").Value = 0.3173459790819593
$ws.Range("J7").Value = 0.3173459790819593
$ws.Range("O7").Value = 0.5929883717175313
$ws.Range("P7").Value = 0.5929883717175314
$ws.Range("Q7").Value = 168.3686026268327
$ws.Range("S7").Value = 0.1881824754069168
$ws.Range("T7").Value = 0.1881824754069168
$ws.Range("I8").Value = 0.2054288400220983
$ws.Range("J8").Value = 0.2054288400220983
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.539707
$ws.Range("N8").Value = 7.619121
$ws.Range("O8").Value = 0.04480768326120512
$ws.Range("P8").Value = 0.04480768326120513
$ws.Range("Q8").Value = 8.235611170318
$ws.Range("R8").Value = 74.120500532862
$ws.Range("S8").Value = 0.009204790396426958
$ws.Range("T8").Value = 0.00920479039642696
$ws.Range("I9").Value = 0.2054288400220983
$ws.Range("J9").Value = 0.2054288400220983
$ws.Range("O9").Value = 0.3622039450212636
$ws.Range("P9").Value = 0.3622039450212636
$ws.Range("S9").Value = 0.07440713627714604
$ws.Range("T9").Value = 0.07440713627714604
$ws.Range("I10").Value = 0.2054288400220983
$ws.Range("J10").Value = 0.2054288400220983
$ws.Range("O10").Value = 0.5929883717175313
$ws.Range("P10").Value = 0.5929883717175314
$ws.Range("S10").Value = 0.1218169133485253
$ws.Range("T10").Value = 0.1218169133485253
